# Update automàtic: dades i banners [2026-02-05 18:19]
# Applies the scraped Meteocat data refresh to Dades_Meteo sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value without Excel re-interpreting it
# as a number/percentage (e.g. "73%" -> 0.73). We build the text via a
# literal-string formula, then Copy + PasteSpecial(Values) to collapse it
# back down to a plain static value while preserving the cell style.
function Set-LiteralText {
    param($addr, $value)
    $escaped = $value -replace '"', '""'
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}
$excel.CutCopyMode = 0

$ws.Range("E2").Value = "2026-02-05 18:17:56"
$ws.Range("E3").Value = "2026-02-05 18:17:59"
$ws.Range("I3").Value = "9.1 mm"
$ws.Range("E4").Value = "2026-02-05 18:18:01"
$ws.Range("J4").Value = "990.2 hPa"
$ws.Range("O4").Value = "11.1 °C"
$ws.Range("E5").Value = "2026-02-05 18:18:04"
Set-LiteralText "H5" '73%'
$ws.Range("J5").Value = "990.3 hPa"
$ws.Range("O5").Value = "9.4 °C"
$ws.Range("E6").Value = "2026-02-05 18:18:06"
$ws.Range("J6").Value = "991.9 hPa"
$ws.Range("M6").Value = "16.0 °C 17:58 TU"
$ws.Range("O6").Value = "12.6 °C"
$ws.Range("E7").Value = "2026-02-05 18:18:09"
Set-LiteralText "H7" '81%'
$ws.Range("J7").Value = "991.8 hPa"
$ws.Range("E8").Value = "2026-02-05 18:18:11"
Set-LiteralText "H8" '88%'
$ws.Range("O8").Value = "8.3 °C"
$ws.Range("E9").Value = "2026-02-05 18:18:14"
$ws.Range("O9").Value = "2.1 °C"
$ws.Range("E10").Value = "2026-02-05 18:18:16"
Set-LiteralText "H10" '89%'
$ws.Range("O10").Value = "7.5 °C"
$ws.Range("E11").Value = "2026-02-05 18:18:18"
$ws.Range("J11").Value = "995.2 hPa"
$ws.Range("M11").Value = "3.6 °C 17:54 TU"
$ws.Range("E12").Value = "2026-02-05 18:18:21"
Set-LiteralText "H12" '91%'
$ws.Range("L12").Value = "42.5 km/h - 264º 17:49 TU"
$ws.Range("M12").Value = "16.3 °C 17:51 TU"
$ws.Range("O12").Value = "9.5 °C"
$ws.Range("E13").Value = "2026-02-05 18:18:24"
$ws.Range("E14").Value = "2026-02-05 18:18:26"
Set-LiteralText "H14" '71%'
$ws.Range("I14").Value = "6.5 mm"
$ws.Range("E15").Value = "2026-02-05 18:18:28"
$ws.Range("J15").Value = "990.8 hPa"
$ws.Range("O15").Value = "7.6 °C"
$ws.Range("E16").Value = "2026-02-05 18:18:31"
$ws.Range("O16").Value = "3.6 °C"
$ws.Range("E17").Value = "2026-02-05 18:18:34"
$ws.Range("I17").Value = "8.3 mm"
$ws.Range("J17").Value = "995.3 hPa"
$ws.Range("E18").Value = "2026-02-05 18:18:36"
$ws.Range("O18").Value = "-4.3 °C"
$ws.Range("E19").Value = "2026-02-05 18:18:39"
$ws.Range("O19").Value = "7.3 °C"
$ws.Range("E20").Value = "2026-02-05 18:18:41"
$ws.Range("O20").Value = "-1.5 °C"
$ws.Range("E21").Value = "2026-02-05 18:18:44"
$ws.Range("J21").Value = "991.1 hPa"
$ws.Range("O21").Value = "5.9 °C"
$ws.Range("E22").Value = "2026-02-05 18:18:46"
$ws.Range("O22").Value = "8.1 °C"
$ws.Range("E23").Value = "2026-02-05 18:18:49"
$ws.Range("J23").Value = "990.2 hPa"
$ws.Range("O23").Value = "8.1 °C"
$ws.Range("E24").Value = "2026-02-05 18:18:52"
$ws.Range("J24").Value = "989.3 hPa"
$ws.Range("E25").Value = "2026-02-05 18:18:54"
Set-LiteralText "H25" '92%'
$ws.Range("J25").Value = "994.4 hPa"
$ws.Range("E26").Value = "2026-02-05 18:18:57"
$ws.Range("O26").Value = "-0.9 °C"
$ws.Range("E27").Value = "2026-02-05 18:19:00"
$ws.Range("J27").Value = "990.5 hPa"
$ws.Range("O27").Value = "8.5 °C"
$ws.Range("E28").Value = "2026-02-05 18:19:02"
$ws.Range("J28").Value = "993.4 hPa"
$ws.Range("O28").Value = "2.2 °C"
$ws.Range("E29").Value = "2026-02-05 18:19:05"
Set-LiteralText "H29" '82%'
$ws.Range("O29").Value = "8.6 °C"
$ws.Range("E30").Value = "2026-02-05 18:19:07"
$ws.Range("O30").Value = "-1.9 °C"
$ws.Range("E31").Value = "2026-02-05 18:19:10"
$ws.Range("J31").Value = "994.5 hPa"
$ws.Range("E32").Value = "2026-02-05 18:19:12"
$ws.Range("I32").Value = "1.2 mm"
$ws.Range("L32").Value = "55.4 km/h - 288º 17:32 TU"
$ws.Range("E33").Value = "2026-02-05 18:19:15"
$ws.Range("O33").Value = "8.7 °C"
$ws.Range("E34").Value = "2026-02-05 18:19:18"
$ws.Range("O34").Value = "3.4 °C"
$ws.Range("E35").Value = "2026-02-05 18:19:20"
$ws.Range("I35").Value = "4.1 mm"
$ws.Range("O35").Value = "-3.0 °C"
$ws.Range("E36").Value = "2026-02-05 18:19:23"

$excel.CutCopyMode = 0
